# Set column H ("Industries") values to 0 for rows 34 through 80
# (dates 4/2/2020 - 5/18/2020), matching the rest of the H column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H34:H80").Value = 0
